$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 1160.6471
$ws.Range("I41").Value = 1278.6154
$ws.Range("J41").Value = 777.25
$ws.Range("K41").Value = 1278.6154
$ws.Range("L41").Value = 777.25
$ws.Range("M41").Value = -838.6153999999999
$ws.Range("N41").Value = -1657.25
$ws.Range("H53").Value = 1797.5385
$ws.Range("I53").Value = 2540.6667
$ws.Range("J53").Value = 125.5
$ws.Range("K53").Value = 2540.6667
$ws.Range("L53").Value = 125.5
$ws.Range("M53").Value = -1903.6667
$ws.Range("N53").Value = -1399.5
$ws.Range("H62").Value = 3438.625
$ws.Range("I62").Value = 4581.8
$ws.Range("J62").Value = 1533.3334
$ws.Range("K62").Value = 4581.8
$ws.Range("L62").Value = 1533.3334
$ws.Range("M62").Value = -3957.8
$ws.Range("N62").Value = -2781.3334
$ws.Range("H64").Value = 18521514
$ws.Range("I64").Value = 27780478
$ws.Range("J64").Value = 3590
$ws.Range("K64").Value = 27780478
$ws.Range("L64").Value = 3590
$ws.Range("M64").Value = -27780230
$ws.Range("N64").Value = -4086
$ws.Range("H65").Value = 3438.625
$ws.Range("I65").Value = 4581.8
$ws.Range("J65").Value = 1533.3334
$ws.Range("K65").Value = 22909
$ws.Range("L65").Value = 7666.666999999999
$ws.Range("M65").Value = -19789
$ws.Range("N65").Value = -13906.667
$ws.Range("H67").Value = 18521514
$ws.Range("I67").Value = 27780478
$ws.Range("J67").Value = 3590
$ws.Range("K67").Value = 27780478
$ws.Range("L67").Value = 3590
$ws.Range("M67").Value = -27779620
$ws.Range("N67").Value = -5306
$ws.Range("H116").Value = 5769.4443
$ws.Range("I116").Value = 5941.1763
$ws.Range("J116").Value = 2850
$ws.Range("K116").Value = 5941.1763
$ws.Range("L116").Value = 2850
$ws.Range("M116").Value = -2499.1763
$ws.Range("N116").Value = -9734
$ws.Range("H129").Value = 2778.279
$ws.Range("I129").Value = 487.81818
$ws.Range("J129").Value = 3565.625
$ws.Range("K129").Value = 1463.45454
$ws.Range("L129").Value = 10696.875
$ws.Range("M129").Value = 3536.54546
$ws.Range("N129").Value = -20696.875
$ws.Range("H137").Value = 11768136
$ws.Range("I137").Value = 3437.0908
$ws.Range("J137").Value = 33336750
$ws.Range("K137").Value = 10311.2724
$ws.Range("L137").Value = 100010250
$ws.Range("M137").Value = -7761.2724
$ws.Range("N137").Value = -100015350
$ws.Range("H138").Value = 6758685.5
$ws.Range("I138").Value = 1516.579
$ws.Range("J138").Value = 13891253
$ws.Range("K138").Value = 4549.737
$ws.Range("L138").Value = 41673759
$ws.Range("M138").Value = 590.2629999999999
$ws.Range("N138").Value = -41684039

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 59800
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 59800
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 59800
$ws.Range("N62").Value = -61048
$ws.Range("H65").Value = 59800
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 59800
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 179400
$ws.Range("N65").Value = -185640
$ws.Range("H81").Value = 22000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 22000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 22000
$ws.Range("N81").Value = -23996
$ws.Range("H84").Value = 22000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 22000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 66000
$ws.Range("N84").Value = -75984

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9264815
$ws.Range("I31").Value = 5079.6772
$ws.Range("J31").Value = 66675172
$ws.Range("K31").Value = 5079.6772
$ws.Range("L31").Value = 66675172
$ws.Range("M31").Value = -4784.6772
$ws.Range("N31").Value = -66675762
$ws.Range("H34").Value = 9264815
$ws.Range("I34").Value = 5079.6772
$ws.Range("J34").Value = 66675172
$ws.Range("K34").Value = 5079.6772
$ws.Range("L34").Value = 66675172
$ws.Range("M34").Value = -4877.6772
$ws.Range("N34").Value = -66675576
$ws.Range("H62").Value = 2300
$ws.Range("I62").Value = 2300
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 2300
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -1676
$ws.Range("H65").Value = 2300
$ws.Range("I65").Value = 2300
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 11500
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -8380

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 4571.4287
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 4571.4287
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 13714.2861
$ws.Range("N80").Value = -15586.2861
$ws.Range("H83").Value = 4571.4287
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 4571.4287
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 41142.85830000001
$ws.Range("N83").Value = -50502.85830000001
$ws.Range("H122").Value = 2362.3333
$ws.Range("I122").Value = 3098.5
$ws.Range("J122").Value = 890
$ws.Range("K122").Value = 27886.5
$ws.Range("L122").Value = 8010
$ws.Range("M122").Value = -25436.5
$ws.Range("N122").Value = -12910
$ws.Range("H133").Value = 117651256
$ws.Range("I133").Value = 181821400
$ws.Range("J133").Value = 5966.6665
$ws.Range("K133").Value = 545464200
$ws.Range("L133").Value = 17899.9995
$ws.Range("M133").Value = -545459140
$ws.Range("N133").Value = -28019.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 5804.6665
$ws.Range("I43").Value = 1609.3334
$ws.Range("H107").Value = 2572.25
$ws.Range("I107").Value = 3662.875
$ws.Range("J107").Value = 1481.625
$ws.Range("K107").Value = 3662.875
$ws.Range("L107").Value = 1481.625
$ws.Range("M107").Value = -1742.875
$ws.Range("N107").Value = -5321.625
$ws.Range("H113").Value = 1590
$ws.Range("I113").Value = 1469.1428
$ws.Range("J113").Value = 2013
$ws.Range("K113").Value = 1469.1428
$ws.Range("L113").Value = 2013
$ws.Range("M113").Value = 700.8571999999999
$ws.Range("N113").Value = -6353

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7131.231
$ws.Range("I40").Value = 13825.25
$ws.Range("J40").Value = 4156.1113
$ws.Range("K40").Value = 13825.25
$ws.Range("L40").Value = 4156.1113
$ws.Range("M40").Value = -13689.25
$ws.Range("N40").Value = -4428.1113
$ws.Range("H132").Value = 7048205
$ws.Range("I132").Value = 3942.4038
$ws.Range("J132").Value = 26327240
$ws.Range("K132").Value = 11827.2114
$ws.Range("L132").Value = 78981720
$ws.Range("M132").Value = -9297.2114
$ws.Range("N132").Value = -78986780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2232.182
$ws.Range("I113").Value = 342.33334
$ws.Range("J113").Value = 4500
$ws.Range("K113").Value = 1027.00002
$ws.Range("L113").Value = 13500
$ws.Range("M113").Value = 1142.99998
$ws.Range("N113").Value = -17840
$ws.Range("H126").Value = 3963.1667
$ws.Range("I126").Value = 2022.4667
$ws.Range("J126").Value = 13666.667
$ws.Range("K126").Value = 6067.4001
$ws.Range("L126").Value = 41000.001
$ws.Range("M126").Value = -3597.4001
$ws.Range("N126").Value = -45940.001
$ws.Range("H132").Value = 1685.2122
$ws.Range("I132").Value = 1271.24
$ws.Range("J132").Value = 2978.875
$ws.Range("K132").Value = 3813.72
$ws.Range("L132").Value = 8936.625
$ws.Range("M132").Value = -1283.72
$ws.Range("N132").Value = -13996.625
